$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "Samtec"
$ws.Range("E3").Value = "SSW-115-01-T-S"
$ws.Range("F3").Value = "CONN RCPT 15POS 0.1 TIN PCB"

$ws.Range("F22").Select()
